$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "26.493.86"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").Value = "1.808.99"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.004"
$ws.Range("E5").Value = "  -0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "307.26"
$ws.Range("E6").Value = "  -1.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4528"
$ws.Range("E7").Value = "  -1.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3602"
$ws.Range("E8").Value = "  -1.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.53"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07096"
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8902"
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07831"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Value = "1.825.13"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.337"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "85.24"
$ws.Range("E17").Value = "  -2.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.006"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008513"
$ws.Range("E19").Value = "  -2.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "26.531.88"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.28"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.980"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "2.051.30"
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.54"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.972"
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.93"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.85"
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.057"
$ws.Range("E29").Value = "  +3.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "112.15"
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.885"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08711"
$ws.Range("E32").Value = "  -1.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.121"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.905"
$ws.Range("E34").Value = "  +15.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.447"
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7234"
$ws.Range("E36").Value = "  -2.35%  "
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.003"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01939"
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5175"
$ws.Range("E43").Value = "  +4.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.789"
$ws.Range("E44").Value = "  -1.56%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1514"
$ws.Range("E45").Value = "  -4.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.031"
$ws.Range("E46").Value = "  -2.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4673"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  -0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.956"
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.04"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.575"
$ws.Range("E51").Value = "  -1.98%  "
